# Session 5 ("Model Governance") time slot moved from 14:30 to 14:50,
# and the workshop date placeholders were refreshed from
# "Monday, June 29, 2020" to "Wednesday, July 1, 2020".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the "14:30 - 16:00" time range inside the Agenda (afternoon)
#    SmartArt graphic to "14:50 - 16:00".
# ---------------------------------------------------------------------
$dash = [char]0x2013
$oldRange = "14:30 " + $dash + " 16:00"
$newRange = "14:50 " + $dash + " 16:00"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasSmartArt) {
            $allNodes = $shape.SmartArt.AllNodes
            for ($ni = 1; $ni -le $allNodes.Count; $ni++) {
                $node = $allNodes.Item($ni)
                $tr = $node.TextFrame2.TextRange
                if ($tr.Text -eq $oldRange) {
                    $tr.Text = $newRange
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the cached "Date Placeholder" text (the datetime2 field)
#    on the slide master and every slide layout.
# ---------------------------------------------------------------------
$oldDate = "Monday, June 29, 2020"
$newDate = "Wednesday, July 1, 2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        # Every "Date Placeholder" shape in this deck caches the same
        # fixed-format field text ("Monday, June 29, 2020"); refresh it
        # unconditionally rather than gating on an exact text match
        # (one of the shapes round-trips its current text as empty
        # through TextRange.Text, even though it holds the same field).
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}
